$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Create "HistDisplayValues" as a copy of "DisplayValues", placed right
# after "SignificanceValues" (i.e. before "Lookups").
# ---------------------------------------------------------------------------
$dispSrc = $wb.Worksheets.Item("DisplayValues")
$dispSrc.Copy($null, $wb.Worksheets.Item("SignificanceValues"))
$histDisp = $wb.Worksheets.Item(4)
$histDisp.Name = "HistDisplayValues"

# ---------------------------------------------------------------------------
# Create "HistSignificanceValues" as a copy of "SignificanceValues", placed
# right after "HistDisplayValues" (i.e. still before "Lookups").
# ---------------------------------------------------------------------------
$sigSrc = $wb.Worksheets.Item("SignificanceValues")
$sigSrc.Copy($null, $wb.Worksheets.Item("HistDisplayValues"))
$histSig = $wb.Worksheets.Item(5)
$histSig.Name = "HistSignificanceValues"

# ---------------------------------------------------------------------------
# Add the defined names that expose the new historical sheets, mirroring
# the ones that already exist for DisplayValues / SignificanceValues.
# (Backtick-escape the "$" so PowerShell doesn't try to interpolate $B/$C
# etc. as variables inside the double-quoted strings.)
# ---------------------------------------------------------------------------
$wb.Names.Add('hist_disp_value_col_head', "='HistDisplayValues'!`$B`$1:`$C`$1")
$wb.Names.Add('hist_disp_value_row_head', "='HistDisplayValues'!`$A`$2:`$A`$3")
$wb.Names.Add('hist_disp_value_values',   "='HistDisplayValues'!`$B`$2:`$C`$3")
$wb.Names.Add('hist_sig_value_col_head',  "='HistSignificanceValues'!`$B`$1:`$C`$1")
$wb.Names.Add('hist_sig_value_row_head',  "='HistSignificanceValues'!`$A`$2:`$A`$3")
$wb.Names.Add('hist_sig_value_values',    "='HistSignificanceValues'!`$B`$2:`$C`$3")
